$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Helper: type a sequence of text "parts" into a cell so that each part
# becomes its own <w:r> run (mirroring how the original authors' edits
# produced multiple adjacent runs instead of one merged run). This is done
# by briefly enabling track-changes while typing (each TypeText call then
# becomes a distinct tracked insertion / run) and then accepting the
# individual revisions one by one, which keeps the run boundaries intact.
function Insert-SplitText($cell, [string[]]$parts) {
    $rng = $cell.Range
    $rng.Select()
    $sel = $word.Selection
    $sel.Collapse(1)

    $d.TrackRevisions = $true
    foreach ($p in $parts) {
        $sel.TypeText($p)
    }
    $d.TrackRevisions = $false

    for ($i = $d.Revisions.Count; $i -ge 1; $i--) {
        $d.Revisions.Item($i).Accept()
    }
}

# --- "Deckard" / "RQ1" row ---
# Ctags column: empty -> "10 m 41 s" split across 6 runs
Insert-SplitText $t.Cell(6, 3) @("1", "0", " m", " ", "41", " s")
# Carol column: empty -> "50h 10m 30s"
Insert-SplitText $t.Cell(6, 6) @("50h 10m 30s")
# Jabref column: empty -> "1 d + "
Insert-SplitText $t.Cell(6, 7) @("1 d + ")

# --- "Deckard" / "RQ2" row ---
# Ctags column: empty -> "19 m"
Insert-SplitText $t.Cell(7, 3) @("19 m")

# BrlCad column: "38 m" -> "36 m" split across 3 runs ("3", "6", " m")
$brlcad = $t.Cell(7, 4)
$brlcad.Range.Text = ""
Insert-SplitText $brlcad @("3", "6", " m")
